$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new rows for new chapters "12.8" (before 15.1), "16.1" (before 16.2),
# and "19.3" (before 19.4). Row numbers are given in document order as each
# insert shifts subsequent rows down by one.
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(25).Insert()
$ws.Rows.Item(29).Insert()

# Set the chapter labels for the newly inserted rows
$ws.Range("A22").Value = "12.8"
$ws.Range("A25").Value = "16.1"
$ws.Range("A29").Value = "19.3"

# Add the "link" markers (value 1) for the chapters that now have links
$ws.Range("G22").Value = 1
$ws.Range("E25").Value = 1
$ws.Range("E26").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("D29").Value = 1

# Re-apply the shared "has any link" formula across the whole data range as a
# single operation (so Excel keeps it as one shared formula group), then clear
# the two rows (25 and 29) that intentionally have no N formula.
$ws.Range("N4:N37").Formula = "=IF(SUM(B4:M4)>0,1,0)"
$ws.Range("N25").ClearContents()
$ws.Range("N29").ClearContents()

# Update the summary formula in O1 to reflect the new data range
$ws.Range("O1").Formula = "=SUM(N3:N1000)/COUNT(N3:N1000)"

# Update conditional formatting range to cover the new rows
$cf = $ws.Range("N3:N34").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("N3:N37"))

# Update the selected cell
$ws.Range("I24").Select() | Out-Null
